# Auto-generated edit script applying value updates to match target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19770

$ws.Range("H33").Value = 377.35
$ws.Range("I33").Value = 308.33334
$ws.Range("K33").Value = 308.33334
$ws.Range("M33").Value = -79.33334000000002

$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19621

$ws.Range("H43").Value = 1250
$ws.Range("J43").Value = 1500
$ws.Range("L43").Value = 1500
$ws.Range("N43").Value = -1638

$ws.Range("H112").Value = 2699.5625
$ws.Range("J112").Value = 2699.5625
$ws.Range("L112").Value = 8098.6875
$ws.Range("N112").Value = -10314.6875

$ws.Range("H132").Value = 22851.137
$ws.Range("I132").Value = 3591.3555
$ws.Range("J132").Value = 167299.5
$ws.Range("K132").Value = 10774.0665
$ws.Range("L132").Value = 501898.5
$ws.Range("M132").Value = -8244.066500000001
$ws.Range("N132").Value = -506958.5

$ws.Range("H137").Value = 10794.151
$ws.Range("I137").Value = 2905.2222
$ws.Range("J137").Value = 13752.5
$ws.Range("K137").Value = 8715.6666
$ws.Range("L137").Value = 41257.5
$ws.Range("M137").Value = -6165.6666
$ws.Range("N137").Value = -46357.5

$ws.Range("H138").Value = 5871.763
$ws.Range("J138").Value = 5366
$ws.Range("L138").Value = 16098
$ws.Range("N138").Value = -26378

$ws.Range("H141").Value = 3166.0833
$ws.Range("I141").Value = 3076.4443
$ws.Range("J141").Value = 3435
$ws.Range("K141").Value = 9229.332900000001
$ws.Range("L141").Value = 10305
$ws.Range("M141").Value = -4049.332900000001
$ws.Range("N141").Value = -20665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1436473.1
$ws.Range("I32").Value = 640221.9
$ws.Range("K32").Value = 640221.9
$ws.Range("M32").Value = -639934.9

$ws.Range("H74").Value = 34484840
$ws.Range("I74").Value = 1692.875
$ws.Range("K74").Value = 1692.875
$ws.Range("M74").Value = -818.875

$ws.Range("H77").Value = 34484840
$ws.Range("I77").Value = 1692.875
$ws.Range("K77").Value = 8464.375
$ws.Range("M77").Value = -4096.375

$ws.Range("H132").Value = 2957.8408
$ws.Range("I132").Value = 2319.9644
$ws.Range("J132").Value = 4074.125
$ws.Range("K132").Value = 6959.8932
$ws.Range("L132").Value = 12222.375
$ws.Range("M132").Value = -4429.8932
$ws.Range("N132").Value = -17282.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2608331
$ws.Range("I31").Value = 4089.4285
$ws.Range("J31").Value = 3052957.5
$ws.Range("K31").Value = 4089.4285
$ws.Range("L31").Value = 3052957.5
$ws.Range("M31").Value = -3794.4285
$ws.Range("N31").Value = -3053547.5

$ws.Range("H34").Value = 2608331
$ws.Range("I34").Value = 4089.4285
$ws.Range("J34").Value = 3052957.5
$ws.Range("K34").Value = 4089.4285
$ws.Range("L34").Value = 3052957.5
$ws.Range("M34").Value = -3887.4285
$ws.Range("N34").Value = -3053361.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6500
$ws.Range("J75").Value = 6500
$ws.Range("L75").Value = 19500
$ws.Range("N75").Value = -21496

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H78").Value = 6500
$ws.Range("J78").Value = 6500
$ws.Range("L78").Value = 58500
$ws.Range("N78").Value = -68484

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H107").Value = 2989.5557
$ws.Range("J107").Value = 3238.375
$ws.Range("L107").Value = 9715.125
$ws.Range("N107").Value = -13555.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27500.5
$ws.Range("J46").Value = 50001
$ws.Range("L46").Value = 50001
$ws.Range("N46").Value = -50313

$ws.Range("H102").Value = 30310298
$ws.Range("I102").Value = 142857780
$ws.Range("J102").Value = 9053.308000000001
$ws.Range("K102").Value = 142857780
$ws.Range("L102").Value = 9053.308000000001
$ws.Range("M102").Value = -142856158
$ws.Range("N102").Value = -12297.308

$ws.Range("H122").Value = 27030268
$ws.Range("I122").Value = 2865.5
$ws.Range("J122").Value = 90913224
$ws.Range("K122").Value = 8596.5
$ws.Range("L122").Value = 272739672
$ws.Range("M122").Value = -6146.5
$ws.Range("N122").Value = -272744572

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 1406.4445
$ws.Range("I132").Value = 1461.6154
$ws.Range("J132").Value = 1263
$ws.Range("K132").Value = 4384.8462
$ws.Range("L132").Value = 3789
$ws.Range("M132").Value = -1854.8462
$ws.Range("N132").Value = -8849

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2233.3333
$ws.Range("I22").Value = 2233.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2233.3333
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 2233.3333
$ws.Range("I27").Value = 2233.3333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2233.3333
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H82").Value = 1470.5
$ws.Range("I82").Value = 1356.8572
$ws.Range("J82").Value = 1629.6
$ws.Range("K82").Value = 1356.8572
$ws.Range("L82").Value = 1629.6
$ws.Range("M82").Value = -995.8571999999999
$ws.Range("N82").Value = -2351.6

$ws.Range("H85").Value = 1470.5
$ws.Range("I85").Value = 1356.8572
$ws.Range("J85").Value = 1629.6
$ws.Range("K85").Value = 1356.8572
$ws.Range("L85").Value = 1629.6
$ws.Range("M85").Value = -108.8571999999999
$ws.Range("N85").Value = -4125.6

$ws.Range("H107").Value = 5495.5
$ws.Range("I107").Value = 5495.5
$ws.Range("K107").Value = 5495.5
$ws.Range("M107").Value = -3575.5

$ws.Range("H132").Value = 10755.125
$ws.Range("I132").Value = 7153
$ws.Range("J132").Value = 15386.429
$ws.Range("K132").Value = 21459
$ws.Range("L132").Value = 46159.287
$ws.Range("M132").Value = -18929
$ws.Range("N132").Value = -51219.287

$ws.Range("H135").Value = 86333.336
$ws.Range("J135").Value = 86333.336
$ws.Range("L135").Value = 86333.336
$ws.Range("N135").Value = -96473.336

$ws.Range("H136").Value = 6981.2
$ws.Range("I136").Value = 5801.1665
$ws.Range("K136").Value = 17403.4995
$ws.Range("M136").Value = -14853.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H81").Value = 7249.95
$ws.Range("I81").Value = 14333.333
$ws.Range("J81").Value = 5999.9414
$ws.Range("K81").Value = 28666.666
$ws.Range("L81").Value = 11999.8828
$ws.Range("M81").Value = -27605.666
$ws.Range("N81").Value = -14121.8828

$ws.Range("H84").Value = 7249.95
$ws.Range("I84").Value = 14333.333
$ws.Range("J84").Value = 5999.9414
$ws.Range("K84").Value = 143333.33
$ws.Range("L84").Value = 59999.414
$ws.Range("M84").Value = -138029.33
$ws.Range("N84").Value = -70607.41399999999

$ws.Range("H132").Value = 1952.1305
$ws.Range("I132").Value = 1999.95
$ws.Range("J132").Value = 1633.3334
$ws.Range("K132").Value = 5999.85
$ws.Range("L132").Value = 4900.0002
$ws.Range("M132").Value = -3469.85
$ws.Range("N132").Value = -9960.0002

$ws.Range("H136").Value = 8337203.5
$ws.Range("I136").Value = 9806693
$ws.Range("J136").Value = 10097.667
$ws.Range("K136").Value = 29420079
$ws.Range("L136").Value = 30293.001
$ws.Range("M136").Value = -29417529
$ws.Range("N136").Value = -35393.001
